$d = $word.ActiveDocument

$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertXML('<w:p><w:pPr><w:rPr><w:lang w:val="fr-CA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>Le sujet du TP2 2 est le même que celui du TP1. L''</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>entiéreté</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t xml:space="preserve"> du projet a été modifié en MVC.</w:t></w:r></w:p>')

$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertXML('<w:p><w:pPr><w:rPr><w:lang w:val="fr-CA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t xml:space="preserve">Les catégories et </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>tages</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t xml:space="preserve"> peuvent être modifiés.</w:t></w:r></w:p>')

$p8 = $d.Paragraphs.Item(8)
$p8.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="Heading2"/><w:rPr><w:lang w:val="fr-CA"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>WebDev</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>')

$p10 = $d.Paragraphs.Item(10)
$p10.Range.InsertXML('<w:p><w:pPr><w:rPr><w:lang w:val="fr-CA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t xml:space="preserve">J’ai tenté de mettre des </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>echos</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t xml:space="preserve"> afin de déterminer si</w:t></w:r></w:p>')

$p11 = $d.Paragraphs.Item(11)
$p11.Range.InsertXML('<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="fr-CA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>On entre dans /route/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>web.php</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t> : OUI</w:t></w:r></w:p>')

$p14 = $d.Paragraphs.Item(14)
$p14.Range.InsertXML('<w:p><w:pPr><w:rPr><w:lang w:val="fr-CA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>Seul la route ‘’ (sans url) retourne une page. Si un autre contrôleur est mis dans la route ‘’, la page est retournée.</w:t></w:r></w:p>')

$p15 = $d.Paragraphs.Item(15)
$p15.Range.InsertXML('<w:p><w:pPr><w:rPr><w:lang w:val="fr-CA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>Le message d’erreur indique que le fichier ne peut être retourné.</w:t></w:r></w:p>')

$p16 = $d.Paragraphs.Item(16)
$p16.Range.InsertXML('<w:p><w:pPr><w:rPr><w:lang w:val="fr-CA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>L''erreur semble venir du fichier .</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>htaccess</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t xml:space="preserve">. Sur </w:t></w:r><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>la version locale</w:t></w:r><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t xml:space="preserve">, si l''utilisateur </w:t></w:r><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>écrit</w:t></w:r><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t xml:space="preserve"> n''importe quoi après le chemin vers le dossier contenant le </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>index.php</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t xml:space="preserve">, une page erreur 404 est retournée. Cependant, sur </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>webdev</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t xml:space="preserve">, le message est "File not </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>found</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>.".</w:t></w:r></w:p>')

$p18 = $d.Paragraphs.Item(18)
$p18.Range.InsertXML('<w:p><w:pPr><w:rPr><w:lang w:val="fr-CA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t xml:space="preserve">Le dossier Diagram contient le </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>schema</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t xml:space="preserve"> ER, et les requêtes SQL pour créer les tables dans </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>workbench</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t xml:space="preserve"> et dans phpMyAdmin (_</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>webdev</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>).</w:t></w:r></w:p>')

$p20 = $d.Paragraphs.Item(20)
$p20.Range.InsertXML('<w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Webdev</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>: https://e2396414.webdev.cmaisonneuve.qc.ca/blog/</w:t></w:r></w:p>')

$p26 = $d.Paragraphs.Item(26)
$p26.Range.InsertXML('<w:p><w:pPr><w:rPr><w:lang w:val="fr-CA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t xml:space="preserve">Chaque article peut seulement avoir été écrit par une personne. Ils contiennent </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>minalement</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t xml:space="preserve"> un titre, un contenu. Les tags et les catégories sont facultatifs. Les catégories sont tirées d''une liste établie par les administrateurs. Les tags peuvent être ajoutés librement. La liste est "case-sensitive". Les autres champs tels la date de création sont remplis automatiquement à la création de l''article.</w:t></w:r></w:p>')

$p27 = $d.Paragraphs.Item(27)
$p27.Range.InsertXML('<w:p><w:pPr><w:rPr><w:lang w:val="fr-CA"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>Le titre, contenu d''un article peuvent être modifié par l''auteur de l''article. L''</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>horodage</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t xml:space="preserve"> de mise à jour est modifié. L''</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t>horodage</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="fr-CA"/></w:rPr><w:t xml:space="preserve"> de la création reste inchangée.</w:t></w:r></w:p>')

Write-Output "done"